$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 403.32
$ws.Range("I92").Value = 225.76471
$ws.Range("J92").Value = 780.625
$ws.Range("K92").Value = 225.76471
$ws.Range("L92").Value = 780.625
$ws.Range("M92").Value = 1022.23529
$ws.Range("N92").Value = -3276.625

# Row 129: Practical Command
$ws.Range("H129").Value = 2364
$ws.Range("I129").Value = 3012.125
$ws.Range("J129").Value = 1499.8334
$ws.Range("K129").Value = 9036.375
$ws.Range("L129").Value = 4499.5002
$ws.Range("M129").Value = -4036.375
$ws.Range("N129").Value = -14499.5002

# Row 138: All-night Crafting
$ws.Range("H138").Value = 3128.923
$ws.Range("I138").Value = 1532.8948
$ws.Range("J138").Value = 3788.152
$ws.Range("K138").Value = 4598.6844
$ws.Range("L138").Value = 11364.456
$ws.Range("M138").Value = 541.3155999999999
$ws.Range("N138").Value = -21644.456

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2338.16
$ws.Range("I61").Value = 2072.7
$ws.Range("J61").Value = 3400
$ws.Range("K61").Value = 2072.7
$ws.Range("L61").Value = 3400
$ws.Range("M61").Value = -1860.7
$ws.Range("N61").Value = -3824

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 10001233
$ws.Range("I74").Value = 12196150
$ws.Range("J74").Value = 2167.6667
$ws.Range("K74").Value = 12196150
$ws.Range("L74").Value = 2167.6667
$ws.Range("M74").Value = -12195276
$ws.Range("N74").Value = -3915.6667

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 10001233
$ws.Range("I77").Value = 12196150
$ws.Range("J77").Value = 2167.6667
$ws.Range("K77").Value = 60980750
$ws.Range("L77").Value = 10838.3335
$ws.Range("M77").Value = -60976382
$ws.Range("N77").Value = -19574.3335

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2338.16
$ws.Range("I136").Value = 2072.7
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 6218.099999999999
$ws.Range("L136").Value = 10200
$ws.Range("M136").Value = -3668.099999999999
$ws.Range("N136").Value = -15300

# Row 139: Backing up My Words
$ws.Range("H139").Value = 33957.223
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 33957.223
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 33957.223
$ws.Range("N139").Value = -44237.223

$ws = $wb.Worksheets.Item("BSM")
# Row 62: Barring the Gates to Foundation
$ws.Range("H62").Value = 31427.857
$ws.Range("I62").Value = 20000
$ws.Range("J62").Value = 33332.5
$ws.Range("K62").Value = 20000
$ws.Range("L62").Value = 33332.5
$ws.Range("M62").Value = -19314
$ws.Range("N62").Value = -34704.5

# Row 65: Starting Young (L)
$ws.Range("H65").Value = 31427.857
$ws.Range("I65").Value = 20000
$ws.Range("J65").Value = 33332.5
$ws.Range("K65").Value = 60000
$ws.Range("L65").Value = 99997.5
$ws.Range("M65").Value = -56568
$ws.Range("N65").Value = -106861.5

# Row 107: The Gold Experience
$ws.Range("H107").Value = 972.9524
$ws.Range("I107").Value = 954.82355
$ws.Range("J107").Value = 1050
$ws.Range("K107").Value = 954.82355
$ws.Range("L107").Value = 1050
$ws.Range("M107").Value = 965.17645
$ws.Range("N107").Value = -4890

# Row 135: Axes to the Maxes
$ws.Range("H135").Value = 61188.89
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 61188.89
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 61188.89
$ws.Range("N135").Value = -71328.89

# Row 138: Bladewinner
$ws.Range("H138").Value = 34628
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 34628
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 34628
$ws.Range("N138").Value = -44908

$ws = $wb.Worksheets.Item("CRP")
# Row 25: Bowing to Necessity
$ws.Range("H25").Value = 30000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 30000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 30000
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -30348

$ws = $wb.Worksheets.Item("CUL")
# Row 103: West Meats East
$ws.Range("H103").Value = 2071.138
$ws.Range("I103").Value = 789
$ws.Range("J103").Value = 2338.25
$ws.Range("K103").Value = 2367
$ws.Range("L103").Value = 7014.75
$ws.Range("M103").Value = -1488
$ws.Range("N103").Value = -8772.75

# Row 122: Salt of the North
$ws.Range("H122").Value = 833.5161000000001
$ws.Range("I122").Value = 819.0833
$ws.Range("J122").Value = 842.6316
$ws.Range("K122").Value = 7371.7497
$ws.Range("L122").Value = 7583.6844
$ws.Range("M122").Value = -4921.7497
$ws.Range("N122").Value = -12483.6844

# Row 129: Comfort Food
$ws.Range("H129").Value = 1283.75
$ws.Range("I129").Value = 472.5
$ws.Range("J129").Value = 2095
$ws.Range("K129").Value = 1417.5
$ws.Range("L129").Value = 6285
$ws.Range("M129").Value = 3582.5
$ws.Range("N129").Value = -16285

# Row 136: Simple Is Hardest
$ws.Range("H136").Value = 792.50507

# Row 137: Creative Chocolate
$ws.Range("H137").Value = 18573640
$ws.Range("I137").Value = 41668084
$ws.Range("J137").Value = 1252806.1
$ws.Range("K137").Value = 125004252
$ws.Range("L137").Value = 3758418.3
$ws.Range("M137").Value = -124999152
$ws.Range("N137").Value = -3768618.3

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar
$ws.Range("H132").Value = 4674.696
$ws.Range("I132").Value = 5673.4814
$ws.Range("J132").Value = 3255.3684
$ws.Range("K132").Value = 17020.4442
$ws.Range("L132").Value = 9766.1052
$ws.Range("M132").Value = -14490.4442
$ws.Range("N132").Value = -14826.1052

# Row 135: Fan of the Foreign
$ws.Range("H135").Value = 39846.668
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 39846.668
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 39846.668
$ws.Range("N135").Value = -49986.668

# Row 137: Sew Excited
$ws.Range("H137").Value = 38036
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 38036
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 38036
$ws.Range("N137").Value = -48236

# Row 138: Orders Anonymous
$ws.Range("H138").Value = 62400
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 62400
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 62400
$ws.Range("N138").Value = -72680

# Row 139: Ringing Gratitude
$ws.Range("H139").Value = 45600
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 45600
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 45600
$ws.Range("N139").Value = -55880

# Row 140: The Right Rod
$ws.Range("H140").Value = 47882.855
$ws.Range("I140").Value = 60000
$ws.Range("J140").Value = 45863.332
$ws.Range("K140").Value = 60000
$ws.Range("L140").Value = 45863.332
$ws.Range("M140").Value = -54820
$ws.Range("N140").Value = -56223.332

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 440.5
$ws.Range("I22").Value = 223.8125
$ws.Range("J22").Value = 729.4167
$ws.Range("K22").Value = 223.8125
$ws.Range("L22").Value = 729.4167
$ws.Range("M22").Value = 71.1875
$ws.Range("N22").Value = -1319.4167

# Row 27: Fire and Hide
$ws.Range("H27").Value = 440.5
$ws.Range("I27").Value = 223.8125
$ws.Range("J27").Value = 729.4167
$ws.Range("K27").Value = 223.8125
$ws.Range("L27").Value = 729.4167
$ws.Range("M27").Value = -116.8125
$ws.Range("N27").Value = -943.4167

# Row 139: Giving Gatherers Their Gear
$ws.Range("H139").Value = 40787.5
$ws.Range("I139").Value = 30650
$ws.Range("J139").Value = 44166.668
$ws.Range("K139").Value = 30650
$ws.Range("L139").Value = 44166.668
$ws.Range("M139").Value = -25510
$ws.Range("N139").Value = -54446.668

# Row 140: Worqor Zormor or Bust
$ws.Range("H140").Value = 53122
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 53122
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 53122
$ws.Range("N140").Value = -63482

# Row 141: Just Generally Freezing
$ws.Range("H141").Value = 57507.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 57507.25
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 57507.25
$ws.Range("N141").Value = -67867.25

$ws = $wb.Worksheets.Item("WVR")
# Row 46: Crunching the Numbers
$ws.Range("H46").Value = 45000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 45000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 45000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -45462

# Row 134: Cloth for Canvas
$ws.Range("H134").Value = 45000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 45000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 135000
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -140070
